{"js": "// Replace the arithmetic-expression text in each cell of the first table,\n// in document order (row-major), with its corresponding new expression.\n// Positional (index-based) mapping is required because some old values\n// (e.g. \"19+18=\") repeat in the grid but map to different replacements\n// depending on which cell they occupy, so cells are addressed by their\n// (row, column) position rather than by a document-wide text search.\nconst replacements = [[\"84-48=\", \"6+0=\"], [\"60-48=\", \"6+45=\"], [\"82-39=\", \"7+77=\"], [\"35-24=\", \"54+2=\"], [\"41-24=\", \"85-35=\"], [\"27+36=\", \"90-13=\"], [\"90-21=\", \"81-54=\"], [\"10+24=\", \"7+19=\"], [\"25+10=\", \"7+61=\"], [\"19+18=\", \"51-51=\"], [\"48+29=\", \"75-52=\"], [\"24+3=\", \"34+39=\"], [\"81+3=\", \"69+13=\"], [\"55+24=\", \"69+0=\"], [\"61-19=\", \"43+3=\"], [\"22-6=\", \"16+33=\"], [\"36-35=\", \"71-32=\"], [\"48+31=\", \"56-12=\"], [\"94-90=\", \"50-27=\"], [\"72-70=\", \"3+62=\"], [\"19+18=\", \"30+45=\"], [\"14+19=\", \"27+11=\"], [\"63-16=\", \"60+19=\"], [\"81-71=\", \"32+50=\"], [\"13+72=\", \"91+2=\"], [\"62+7=\", \"79-8=\"], [\"6+78=\", \"93-44=\"], [\"67+3=\", \"36-36=\"], [\"17+60=\", \"74-65=\"], [\"63+2=\", \"81-44=\"], [\"5+36=\", \"62-43=\"], [\"2+75=\", \"92-41=\"], [\"63+4=\", \"8+41=\"], [\"42-3=\", \"39+53=\"], [\"13+79=\", \"85-42=\"], [\"24+68=\", \"35-4=\"], [\"44-26=\", \"45+25=\"], [\"69+28=\", \"92-89=\"], [\"28+5=\", \"42-24=\"], [\"2+78=\", \"70-42=\"], [\"72-8=\", \"48+28=\"], [\"29+1=\", \"71+21=\"], [\"87-21=\", \"31+66=\"], [\"78-3=\", \"77+9=\"], [\"30+2=\", \"22+64=\"], [\"11+81=\", \"16+50=\"], [\"96-86=\", \"65-10=\"], [\"13+77=\", \"20+20=\"], [\"93-62=\", \"50-24=\"], [\"58-10=\", \"76-37=\"], [\"31+59=\", \"98-59=\"], [\"27+35=\", \"43-29=\"], [\"3+69=\", \"16+72=\"], [\"0+38=\", \"9+78=\"], [\"29+25=\", \"39+22=\"], [\"4+47=\", \"26-15=\"], [\"8+27=\", \"26+47=\"], [\"35-26=\", \"20+14=\"], [\"47-11=\", \"7+56=\"], [\"21+7=\", \"0+30=\"], [\"94-74=\", \"79-34=\"], [\"97-79=\", \"82-77=\"], [\"16+27=\", \"20+73=\"], [\"60-28=\", \"86-23=\"], [\"78-38=\", \"70-58=\"], [\"6+83=\", \"12+8=\"], [\"75-34=\", \"69+17=\"], [\"14+13=\", \"85-59=\"], [\"22+59=\", \"41+47=\"], [\"88-67=\", \"79+18=\"], [\"17+67=\", \"12+17=\"], [\"57-45=\", \"66-0=\"], [\"38+17=\", \"84-45=\"], [\"50-39=\", \"30+66=\"], [\"71-36=\", \"83-10=\"], [\"74-16=\", \"13+57=\"], [\"36-7=\", \"15+81=\"], [\"43-42=\", \"50-34=\"], [\"1+74=\", \"34-12=\"], [\"67-60=\", \"40+49=\"], [\"5+76=\", \"70-21=\"], [\"92-8=\", \"82-78=\"], [\"27+15=\", \"48-8=\"], [\"24+41=\", \"36+42=\"], [\"32+12=\", \"14+54=\"], [\"30+25=\", \"4+28=\"], [\"73-57=\", \"99-90=\"], [\"46+32=\", \"38-31=\"], [\"89-32=\", \"40+16=\"], [\"15+32=\", \"54-28=\"], [\"77-27=\", \"28+52=\"], [\"95-67=\", \"19+62=\"], [\"2+3=\", \"81-6=\"], [\"35+64=\", \"4+72=\"], [\"88-24=\", \"81-48=\"], [\"90-28=\", \"17+35=\"], [\"96-5=\", \"91-40=\"], [\"49+45=\", \"15+42=\"], [\"67-27=\", \"71-67=\"], [\"26+50=\", \"46+53=\"]];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"No tables found in the document body.\");\n}\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\nconst columnCount = 5; // table is a fixed 5-column grid\nconst total = Math.min(table.rowCount * columnCount, replacements.length);\n\n// Pass 1: for every target cell, search (within that cell only) for its\n// expected \"before\" text so the matched range keeps the run's original\n// formatting (font, size, paragraph alignment) when we replace it.\nconst cells = [];\nconst searches = [];\nfor (let idx = 0; idx < total; idx++) {\n  const r = Math.floor(idx / columnCount);\n  const c = idx % columnCount;\n  const [oldText] = replacements[idx];\n\n  const cell = table.getCell(r, c);\n  const found = cell.body.search(oldText, { matchCase: true });\n  found.load(\"items\");\n\n  cells.push(cell);\n  searches.push(found);\n}\nawait context.sync();\n\n// Pass 2: apply each replacement. If the expected old text wasn't found\n// in its cell (unexpected drift), fall back to overwriting the whole\n// cell body so the correct new value still ends up in place.\nfor (let idx = 0; idx < total; idx++) {\n  const [, newText] = replacements[idx];\n  const found = searches[idx];\n\n  if (found.items.length > 0) {\n    found.items[0].insertText(newText, Word.InsertLocation.replace);\n  } else {\n    cells[idx].body.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the arithmetic-expression text in each cell of the first table,\n# in document order (row-major), with its corresponding new expression.\n# Cells are addressed positionally by (row, column) rather than by text\n# search, because some old expressions (e.g. \"19+18=\") repeat in the grid\n# but map to different replacements depending on which cell they occupy.\n\n$oldValues = @(\n    \"84-48=\",\n    \"60-48=\",\n    \"82-39=\",\n    \"35-24=\",\n    \"41-24=\",\n    \"27+36=\",\n    \"90-21=\",\n    \"10+24=\",\n    \"25+10=\",\n    \"19+18=\",\n    \"48+29=\",\n    \"24+3=\",\n    \"81+3=\",\n    \"55+24=\",\n    \"61-19=\",\n    \"22-6=\",\n    \"36-35=\",\n    \"48+31=\",\n    \"94-90=\",\n    \"72-70=\",\n    \"19+18=\",\n    \"14+19=\",\n    \"63-16=\",\n    \"81-71=\",\n    \"13+72=\",\n    \"62+7=\",\n    \"6+78=\",\n    \"67+3=\",\n    \"17+60=\",\n    \"63+2=\",\n    \"5+36=\",\n    \"2+75=\",\n    \"63+4=\",\n    \"42-3=\",\n    \"13+79=\",\n    \"24+68=\",\n    \"44-26=\",\n    \"69+28=\",\n    \"28+5=\",\n    \"2+78=\",\n    \"72-8=\",\n    \"29+1=\",\n    \"87-21=\",\n    \"78-3=\",\n    \"30+2=\",\n    \"11+81=\",\n    \"96-86=\",\n    \"13+77=\",\n    \"93-62=\",\n    \"58-10=\",\n    \"31+59=\",\n    \"27+35=\",\n    \"3+69=\",\n    \"0+38=\",\n    \"29+25=\",\n    \"4+47=\",\n    \"8+27=\",\n    \"35-26=\",\n    \"47-11=\",\n    \"21+7=\",\n    \"94-74=\",\n    \"97-79=\",\n    \"16+27=\",\n    \"60-28=\",\n    \"78-38=\",\n    \"6+83=\",\n    \"75-34=\",\n    \"14+13=\",\n    \"22+59=\",\n    \"88-67=\",\n    \"17+67=\",\n    \"57-45=\",\n    \"38+17=\",\n    \"50-39=\",\n    \"71-36=\",\n    \"74-16=\",\n    \"36-7=\",\n    \"43-42=\",\n    \"1+74=\",\n    \"67-60=\",\n    \"5+76=\",\n    \"92-8=\",\n    \"27+15=\",\n    \"24+41=\",\n    \"32+12=\",\n    \"30+25=\",\n    \"73-57=\",\n    \"46+32=\",\n    \"89-32=\",\n    \"15+32=\",\n    \"77-27=\",\n    \"95-67=\",\n    \"2+3=\",\n    \"35+64=\",\n    \"88-24=\",\n    \"90-28=\",\n    \"96-5=\",\n    \"49+45=\",\n    \"67-27=\",\n    \"26+50=\"\n)\n\n$newValues = @(\n    \"6+0=\",\n    \"6+45=\",\n    \"7+77=\",\n    \"54+2=\",\n    \"85-35=\",\n    \"90-13=\",\n    \"81-54=\",\n    \"7+19=\",\n    \"7+61=\",\n    \"51-51=\",\n    \"75-52=\",\n    \"34+39=\",\n    \"69+13=\",\n    \"69+0=\",\n    \"43+3=\",\n    \"16+33=\",\n    \"71-32=\",\n    \"56-12=\",\n    \"50-27=\",\n    \"3+62=\",\n    \"30+45=\",\n    \"27+11=\",\n    \"60+19=\",\n    \"32+50=\",\n    \"91+2=\",\n    \"79-8=\",\n    \"93-44=\",\n    \"36-36=\",\n    \"74-65=\",\n    \"81-44=\",\n    \"62-43=\",\n    \"92-41=\",\n    \"8+41=\",\n    \"39+53=\",\n    \"85-42=\",\n    \"35-4=\",\n    \"45+25=\",\n    \"92-89=\",\n    \"42-24=\",\n    \"70-42=\",\n    \"48+28=\",\n    \"71+21=\",\n    \"31+66=\",\n    \"77+9=\",\n    \"22+64=\",\n    \"16+50=\",\n    \"65-10=\",\n    \"20+20=\",\n    \"50-24=\",\n    \"76-37=\",\n    \"98-59=\",\n    \"43-29=\",\n    \"16+72=\",\n    \"9+78=\",\n    \"39+22=\",\n    \"26-15=\",\n    \"26+47=\",\n    \"20+14=\",\n    \"7+56=\",\n    \"0+30=\",\n    \"79-34=\",\n    \"82-77=\",\n    \"20+73=\",\n    \"86-23=\",\n    \"70-58=\",\n    \"12+8=\",\n    \"69+17=\",\n    \"85-59=\",\n    \"41+47=\",\n    \"79+18=\",\n    \"12+17=\",\n    \"66-0=\",\n    \"84-45=\",\n    \"30+66=\",\n    \"83-10=\",\n    \"13+57=\",\n    \"15+81=\",\n    \"50-34=\",\n    \"34-12=\",\n    \"40+49=\",\n    \"70-21=\",\n    \"82-78=\",\n    \"48-8=\",\n    \"36+42=\",\n    \"14+54=\",\n    \"4+28=\",\n    \"99-90=\",\n    \"38-31=\",\n    \"40+16=\",\n    \"54-28=\",\n    \"28+52=\",\n    \"19+62=\",\n    \"81-6=\",\n    \"4+72=\",\n    \"81-48=\",\n    \"17+35=\",\n    \"91-40=\",\n    \"15+42=\",\n    \"71-67=\",\n    \"46+53=\"\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$columnCount = 5\n$rowCount = $t.Rows.Count\n\n$idx = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $columnCount; $c++) {\n        if ($idx -ge $newValues.Length) { break }\n        $cell = $t.Cell($r, $c)\n        $cellRange = $cell.Range\n\n        # Sanity check (non-fatal): note if the cell's current text isn't\n        # the expected \"before\" value for this position; we still write\n        # the intended new value either way.\n        $currentText = $cellRange.Text.TrimEnd([char]13, [char]7)\n        $expectedOld = $oldValues[$idx]\n        if ($currentText -ne $expectedOld) {\n            Write-Output \"Warning: cell ($r,$c) was '$currentText', expected '$expectedOld'\"\n        }\n\n        # Assigning to Range.Text replaces just the textual content of the\n        # cell and automatically preserves the trailing cell-mark as well\n        # as the run/paragraph formatting (font, size, alignment) already\n        # applied to that cell.\n        $cellRange.Text = $newValues[$idx]\n\n        $idx++\n    }\n}\n"}
